$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M: "IP Address3" header + 10 IP values ---
$ws.Range("M1").Value = "IP Address3"

$ipAddress3 = @(
    "13.233.168.190",
    "13.201.119.75",
    "13.233.85.17",
    "43.205.231.111",
    "3.108.61.207",
    "13.201.3.165",
    "13.201.45.225",
    "13.233.143.125",
    "15.206.125.144",
    "3.110.133.39"
)

for ($i = 0; $i -lt $ipAddress3.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $ipAddress3[$i]
}

# --- Formatting: reuse existing cell styles via copy/paste-special (formats only) ---
# Header M1 should look like the other header cells (A1's style)
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# Data cells in L (previously a slightly different style) and the new M column
# should match the plain bordered style already used by columns H:K (e.g. K2)
$ws.Range("K2").Copy()
$ws.Range("L2:L11").PasteSpecial(-4122)
$ws.Range("M2:M11").PasteSpecial(-4122)

# --- Column width: match new column M to column L ---
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(12).ColumnWidth

# --- Update view: scrolled so column C is leftmost, selection moves to G19 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G19").Select()
